# Adiciona colunas "modelo" e "politica" (politica de preco) na planilha
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere 2 novas colunas antes da coluna C; as colunas antigas C/D/E (full/tipo/link)
# sao deslocadas para E/F/G automaticamente, preservando seus valores e estilos.
$ws.Columns("C:D").Insert()

# Cabecalhos das novas colunas
$ws.Cells.Item(1, 3).Value2 = "modelo"
$ws.Cells.Item(1, 4).Value2 = "politica"

# Dados por linha: modelo, politica, full, tipo, link
$rows = @(
    @{Row=2; Modelo="FONTE 200A LITE"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:227596415#searchVariation=MLB24154371&position=3&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=3; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-longa-distncia-jfa-k600-completo-preto-e-cinza/p/MLB27970249?pdp_filters=seller_id:227596415#searchVariation=MLB27970249&position=4&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=4; Modelo="FONTE 90 BOB"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:227596415#searchVariation=MLB21562641&position=5&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=5; Modelo="FONTE 200A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:227596415#searchVariation=MLB21348561&position=6&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=6; Modelo="FONTE 200 MONO"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-cor-preto/p/MLB24006449?pdp_filters=seller_id:227596415#searchVariation=MLB24006449&position=7&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=7; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-longa-distancia-jfa-k600-preto-com-verde/p/MLB27999036?pdp_filters=seller_id:227596415#searchVariation=MLB27999036&position=8&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=8; Modelo="FONTE 70A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:227596415#searchVariation=MLB21455208&position=9&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=9; Modelo="FONTE 60A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:227596415#searchVariation=MLB21320712&position=10&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=10; Modelo="FONTE 200 BOB"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-200a-bivolt-cor-preto/p/MLB26854417?pdp_filters=seller_id:227596415#searchVariation=MLB26854417&position=11&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=11; Modelo="FONTE 120A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:227596415#searchVariation=MLB21392652&position=12&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=12; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-vermelho/p/MLB34210379?pdp_filters=seller_id:227596415#searchVariation=MLB34210379&position=13&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=13; Modelo="FONTE 200 BOB"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:227596415#searchVariation=MLB24834408&position=14&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=14; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-preto/p/MLB28687615?pdp_filters=seller_id:227596415#searchVariation=MLB28687615&position=15&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=15; Modelo="FONTE 40A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-automotiva-40-amperes-jfa-storm-red-line-cca-sci-smart-cor-preto/p/MLB21621306?pdp_filters=seller_id:227596415#searchVariation=MLB21621306&position=16&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=16; Modelo="FONTE 60A LITE"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/jfa-fonte-carregador-storm-lite-60a-3000-w-preto/p/MLB23456525?pdp_filters=seller_id:227596415#searchVariation=MLB23456525&position=17&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=17; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-40a-lite-storm-slim-bivolt-cor-preto/p/MLB33435981?pdp_filters=seller_id:227596415#searchVariation=MLB33435981&position=20&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=18; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-longa-distncia-jfa-k1200-pretolaranja-1200mt/p/MLB28357019?pdp_filters=seller_id:227596415#searchVariation=MLB28357019&position=24&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=19; Modelo="FONTE 120A LITE"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-cor-preto/p/MLB23998473?pdp_filters=seller_id:227596415#searchVariation=MLB23998473&position=25&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=20; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-remoto-universal-longa-jfa-k1200-preto-c-verde/p/MLB30476096?pdp_filters=seller_id:227596415#searchVariation=MLB30476096&position=26&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=21; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-200a-storm-lite-mono-220v-cor-azul/p/MLB30464905?pdp_filters=seller_id:227596415#searchVariation=MLB30464905&position=18&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=22; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-remoto-longa-distancia-jfa-k600-preto-laranja/p/MLB31403178?pdp_filters=seller_id:227596415#searchVariation=MLB31403178&position=19&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=23; Modelo="FONTE 200A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://www.mercadolivre.com.br/fonte-carregador-automotiva-storm-sci-redline-jfa-200a-slim-cor-preto/p/MLB26091477?pdp_filters=seller_id:227596415#searchVariation=MLB26091477&position=21&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=24; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/kit-2-controles-longa-distancia-jfa-preto-com-verde-k600-m/p/MLB28056168?pdp_filters=seller_id:227596415#searchVariation=MLB28056168&position=22&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=25; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="premium"; Link="https://www.mercadolivre.com.br/controle-de-longa-distancia-jfa-k600-600-metros-com-10-cores-cor/p/MLB28243528?pdp_filters=seller_id:227596415#searchVariation=MLB28243528&position=23&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=26; Modelo="FONTE 40A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-2640386252-fonte-carregador-automotivo-jfa-storm-40a-bivolt-voltamp-_JM#position%3D27%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=27; Modelo="FONTE 200 MONO"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-2642363094-fonte-jfa-200-storm-monovolt-voltamp-para-modulo-automotivo-_JM#position%3D28%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=28; Modelo="Modelo identificado mas fora do range de preco"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-3288375845-fonte-carregador-automotivo-jfa-60a-storm-lite-12v-bivolt-_JM#position%3D29%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=29; Modelo="FONTE 200A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-3348510540-fonte-carregador-automotivo-jfa-storm-200-amperes-sci-bivolt-_JM#position%3D30%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=30; Modelo="FONTE 120 BOB"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-3348316634-fonte-carregador-automotivo-jfa-bob-storm-120a-bivolt-_JM#position%3D31%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=31; Modelo="FONTE 60A LITE"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-3287832201-fonte-carregador-automotivo-jfa-60a-storm-lite-12v-bivolt-_JM#position%3D32%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=32; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-1011507679-controle-som-longa-distancia-jfa-k1200-azul-alcance-1200-mt-_JM#position%3D33%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=33; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-4051614000-fonte-carregador-automotivo-jfa-40a-storm-lite-12v-bivolt-_JM#position%3D34%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=34; Modelo="FONTE 200 BOB"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-3348494984-fonte-carregador-automotivo-jfa-bob-storm-200a-bivolt-_JM#position%3D35%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=35; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://www.mercadolivre.com.br/controle-longa-distancia-jfa-k1200-alcance-de-1200-metros/p/MLB34245679?pdp_filters=seller_id:227596415#searchVariation=MLB34245679&position=1&search_layout=stack&type=product&tracking_id=2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=36; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-3549434493-fonte-carregador-automotivo-jfa-200a-storm-lite-mono-220v-_JM#position%3D36%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=37; Modelo="FONTE 70A LITE"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-3323988125-fonte-automotiva-jfa-storm-lite-70a-bivolt-carregador-_JM#position%3D37%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=38; Modelo="FONTE 200 BOB"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-3348469250-fonte-carregador-automotivo-jfa-bob-storm-200a-bivolt-_JM?searchVariation=177182710708#searchVariation%3D177182710708%26position%3D38%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=39; Modelo="FONTE 120A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-2717116602-fonte-automotiva-jfa-storm-120a-bivolt-com-medidor-cca-e-sci-_JM#position%3D39%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=40; Modelo="FONTE 70A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-2969484411-fonte-automotiva-jfa-storm-70a-bivolt-com-medidor-cca-e-sci-_JM#position%3D40%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=41; Modelo="FONTE 60A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-2640310770-fonte-carregador-automotivo-jfa-storm-60a-bivolt-voltamp-_JM#position%3D41%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=42; Modelo="FONTE 70A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-2642107531-fonte-jfa-70-storm-bivolt-voltamp-para-modulo-automotivo-_JM#position%3D42%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=43; Modelo="FONTE 70A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-2641738515-fonte-automotiva-jfa-storm-70a-bivolt-com-medidor-cca-e-sci-_JM#position%3D43%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=44; Modelo="FONTE 90 BOB"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-3319995819-fonte-carregador-automotivo-jfa-bob-storm-90a-bivolt-_JM#position%3D44%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=45; Modelo="FONTE 200A LITE"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-3322441391-fonte-carregador-automotivo-jfa-200a-storm-lite-12v-bivolt-_JM#position%3D45%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=46; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-1011572142-kit-2-controles-longa-distancia-jfa-k1200-vermelho-200-mt-_JM#position%3D46%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=47; Modelo="FONTE 120A"; Politica="Igual"; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-2642150646-fonte-jfa-120-storm-bivolt-voltamp-para-modulo-automotivo-_JM#position%3D47%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=48; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="premium"; Link="https://produto.mercadolivre.com.br/MLB-1008127611-kit-2-controles-longa-distancia-jfa-preto-com-cinza-k600-m-_JM#position%3D48%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=49; Modelo="Sem Modelo"; Politica=$null; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-1011566977-controle-longa-distancia-jfa-k1200-vermelho-alcance-1200-mt-_JM#position%3D49%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"},
    @{Row=50; Modelo="FONTE 70A"; Politica="Igual"; Full="NA"; Tipo="classico"; Link="https://produto.mercadolivre.com.br/MLB-2640192838-fonte-carregador-automotivo-jfa-storm-70a-bivolt-voltamp-_JM#position%3D50%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D2017e041-f494-43fb-8fc7-1420f881ec70"}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value2 = $r.Modelo
    if ($null -eq $r.Politica) {
        $ws.Cells.Item($r.Row, 4).Value2 = ""
    } else {
        $ws.Cells.Item($r.Row, 4).Value2 = $r.Politica
    }
    $ws.Cells.Item($r.Row, 5).Value2 = $r.Full
    $ws.Cells.Item($r.Row, 6).Value2 = $r.Tipo
    $ws.Cells.Item($r.Row, 7).Value2 = $r.Link
}
